$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column A for rows 2-4
$ws.Range("A2").Value = 701
$ws.Range("A3").Value = 113
$ws.Range("A4").Value = 602

# Update the active selection to A4
$ws.Range("A4").Select()
